# Apr 2019 to May 2019 - Cord Cutting Wizard - Change Log
# Rename the "Comment" column values (column E) to the new wording scheme.
#
# Old wording -> New wording
#   "Network added to Service in May"      -> "Network Added to Base Service"      (when the New Value is "Yes", i.e. added to the core/base service)
#                                           -> "Network Added to Add-On Package"    (otherwise, i.e. added to a named add-on package)
#   "Network removed from Service in May"  -> "Network Removed from Database"
#   "Service Add-On Package Name Changed"  -> "Name of Add-On Package Changed"
#   "New Network added to DB in May"       -> "New Network Added to Database"

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$oldAdded    = "Network added to Service in May"
$oldRemoved  = "Network removed from Service in May"
$oldPkgName  = "Service Add-On Package Name Changed"
$oldNewDb    = "New Network added to DB in May"

$newBaseAdd  = "Network Added to Base Service"
$newAddOnAdd = "Network Added to Add-On Package"
$newRemoved  = "Network Removed from Database"
$newPkgName  = "Name of Add-On Package Changed"
$newNewDb    = "New Network Added to Database"

# Find the last used row on the sheet (data runs from row 2 through the end).
$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row

for ($r = 2; $r -le $lastRow; $r++) {
    $comment = $ws.Cells.Item($r, 5).Text
    if ($comment -eq $oldAdded) {
        $newValue = $ws.Cells.Item($r, 4).Text
        if ($newValue -eq "Yes") {
            $ws.Cells.Item($r, 5).Value = $newBaseAdd
        } else {
            $ws.Cells.Item($r, 5).Value = $newAddOnAdd
        }
    } elseif ($comment -eq $oldRemoved) {
        $ws.Cells.Item($r, 5).Value = $newRemoved
    } elseif ($comment -eq $oldPkgName) {
        $ws.Cells.Item($r, 5).Value = $newPkgName
    } elseif ($comment -eq $oldNewDb) {
        $ws.Cells.Item($r, 5).Value = $newNewDb
    }
}

# Update the active selection to mirror the recorded cursor position.
$ws.Range("E80").Select()
